$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the formatting of row 8 (Assignment 7) onto new row 9 so the
# existing cellXfs entries (styles 3 and 7) get reused instead of creating
# brand-new style records.
$ws.Range("A8:C8").Copy()
$ws.Range("A9:C9").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Row 9 height matches the other "link" rows (34.5pt, custom height).
$ws.Rows.Item(9).RowHeight = 34.5

# The new B9 cell gets its own style (left + vertical centered) distinct
# from B6:B8 (vertical centered only) -> forces a new cellXfs entry (idx 9).
$ws.Range("B9").HorizontalAlignment = -4131

# Populate the new assignment row.
$ws.Range("A9").Value = "Assignment 8"
$ws.Range("B9").Value = "https://github.com/Vasanth30e/Assignment_Phase_4/tree/master/Assignment_8"
$ws.Range("C9").Value = 45189

# Move the selection the way the author's session ended up.
$null = $ws.Range("B15").Select()
